$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Decrement existing weekly "farms_total_count" / "farms_to_examine_count"
#    (and, for the final existing week 202507, also "farms_examined_count" /
#    "farms_examined_positive_count") values by 1 -- these are the revised
#    historic figures in column D.
# ---------------------------------------------------------------------------
$dUpdates = @{
  "D2"   = 11760
  "D3"   = 11505
  "D7"   = 11864
  "D8"   = 11152
  "D12"  = 11972
  "D13"  = 10650
  "D17"  = 12016
  "D18"  = 10099
  "D22"  = 12073
  "D23"  = 9669
  "D27"  = 12115
  "D28"  = 9209
  "D32"  = 12155
  "D33"  = 8801
  "D37"  = 12193
  "D38"  = 8277
  "D42"  = 12228
  "D43"  = 7709
  "D47"  = 12253
  "D48"  = 7161
  "D52"  = 12281
  "D53"  = 6519
  "D57"  = 12304
  "D58"  = 5909
  "D62"  = 12315
  "D63"  = 5743
  "D67"  = 12337
  "D68"  = 5559
  "D72"  = 12353
  "D73"  = 5123
  "D77"  = 12375
  "D78"  = 4645
  "D82"  = 12393
  "D83"  = 4126
  "D87"  = 12420
  "D88"  = 3720
  "D92"  = 12436
  "D93"  = 3382
  "D97"  = 12457
  "D99"  = 9386
  "D100" = 1549
}
foreach ($addr in $dUpdates.Keys) {
  $ws.Range($addr).Value = $dUpdates[$addr]
}

# ---------------------------------------------------------------------------
# 2) Append the new weekly block for YearWeekIso 202508 (LastDayOfWeek
#    2025-02-23, serial 45711) as rows 102-106, one row per Variable.
# ---------------------------------------------------------------------------
$newWeek = @(
  @(202508, "farms_total_count", 12472),
  @(202508, "farms_to_examine_count", 2795),
  @(202508, "farms_examined_count", 9677),
  @(202508, "farms_examined_positive_count", 1548),
  @(202508, "farms_examined_negative_count", 8129)
)

$targetRow = 102
foreach ($entry in $newWeek) {
  $ws.Cells.Item($targetRow, 1).Value = $entry[0]
  $ws.Cells.Item($targetRow, 2).Value = 45711
  $ws.Cells.Item($targetRow, 3).Value = $entry[1]
  $ws.Cells.Item($targetRow, 4).Value = $entry[2]

  # Reuse the existing date number-format (style index carried by the B
  # column's date cells) instead of minting a new style.
  $ws.Range("B101").Copy()
  $ws.Range("B$targetRow").PasteSpecial(-4122)

  $targetRow = $targetRow + 1
}

# ---------------------------------------------------------------------------
# 3) Rows 107-131: leftover date-formatted (but empty) cells in column B,
#    matching what Excel leaves behind after a fill-down of the date
#    formatting with no further data beneath the new week's rows.
# ---------------------------------------------------------------------------
for ($r = 107; $r -le 131; $r++) {
  $ws.Range("B101").Copy()
  $ws.Range("B$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 4) Update the sheet's recorded selection to match the post-edit UI state.
# ---------------------------------------------------------------------------
$ws.Range("F86").Select() | Out-Null
